# Applies the cryptos.xlsx update described by the commit:
# "Updated cryptos list on Fri Jul 26 08:29:19 UTC 2024 with GitHub Actions"
#
# D-column price cells that parse as plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the source
# workbook, where every Price cell is stored as a string, e.g. "1.00",
# "6.75", "0.0675" must stay literal text, not be coerced to 1, 6.75, 0.0675
# as a float/int).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.083.55"
$ws.Range("E2").Value = "  +4.75%  "
$ws.Range("D3").Value = "3.255.19"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'578.66"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "'177.41"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.601"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "3.254.79"
$ws.Range("E9").Value = "  +2.77%  "
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("D11").Value = "'6.75"
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").Value = "3.819.88"
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'28.11"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("D16").Value = "67.089.09"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("E17").Value = "  +3.08%  "
$ws.Range("D18").Value = "3.255.10"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").Value = "'5.86"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'13.42"
$ws.Range("E20").Value = "  +3.00%  "
$ws.Range("D21").Value = "'373.42"
$ws.Range("E21").Value = "  +5.74%  "
$ws.Range("D22").Value = "'7.64"
$ws.Range("E22").Value = "  +6.12%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'71.10"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "'0.511"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").Value = "3.396.99"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +4.40%  "
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("D33").Value = "'22.61"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +5.40%  "
$ws.Range("D36").Value = "'6.81"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("D37").Value = "'167.03"
$ws.Range("E37").Value = "  +7.83%  "
$ws.Range("E38").Value = "  +4.84%  "
$ws.Range("E39").Value = "  +5.44%  "
$ws.Range("E40").Value = "  +9.94%  "
$ws.Range("D41").Value = "'27.09"
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "2.759.48"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("E44").Value = "  +6.78%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'353.40"
$ws.Range("E45").Value = "  +8.15%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.38"
$ws.Range("E46").Value = "  +5.09%  "
$ws.Range("D47").Value = "'25.34"
$ws.Range("E47").Value = "  +6.04%  "
$ws.Range("D48").Value = "'40.54"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("D49").Value = "'0.0675"
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  +0.67%  "
